$wb = $excel.ActiveWorkbook
$botSide  = $wb.Worksheets.Item("Bot-side")
$hostSide = $wb.Worksheets.Item("Host-side")

# --- Bot-side sheet: row 17 used to be an "unused" instruction slot with a
# note in column V. It is now allocated to the new GET_POS instruction, with
# its description/function moved into the normal Description/Function
# columns (T/U) and the old Notes-column remark cleared out. ---
$botSide.Range("A17").Value = "GET_POS"
$botSide.Range("T17").Value = "Gets bot's X, Y, theta, and timestamp (Special command)"
$botSide.Range("U17").Value = "getPos()"
$botSide.Range("V17").ClearContents()

# --- Host-side sheet: add a new response row (row 6) describing the
# RESP_POS response to the new GET_POS instruction. ---
$hostSide.Range("A6").Value = "RESP_POS"
$hostSide.Range("B6").Value = 0
$hostSide.Range("C6").Value = 0
$hostSide.Range("D6").Value = 0
$hostSide.Range("E6").Value = 0
$hostSide.Range("F6").Value = 0
$hostSide.Range("G6").Value = "x"
$hostSide.Range("H6").Value = "x"
$hostSide.Range("I6").Value = "x"
$hostSide.Range("J6").Value = "…"
$hostSide.Range("S6").Value = "This response has bot ID in most significant byte, followed by the raw 32 bit float values for X, Y, theta, and 32 bit unsigned long. Total is 17 bytes"

$hostSide.Range("S6").WrapText = $true
$hostSide.Rows.Item(6).RowHeight = 48
